$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 20 : "Thây đổi tên cuộc hội thoại cũ"
#   - B20 status Missing -> Available
#   - C20 note filled in (was empty)
# ---------------------------------------------------------------------------
$ws.Range("B13").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value2 = "Available"

$ws.Range("C20").Value2 = "Có thể đổi tên cuộc hội thoại ở dashboard và trang lịch sử"

# ---------------------------------------------------------------------------
# Row 23 : "Đổi tên/Đặt tên cho các cuộc hội thoại"
#   - B23 status Dev -> Available
#   - C23 note replaced
# ---------------------------------------------------------------------------
$ws.Range("B13").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value2 = "Available"

$ws.Range("C23").Value2 = "Có thể đổi tên cuộc hội thoại ở dashboard và trang lịch sử"

# ---------------------------------------------------------------------------
# Row 12 : "Cập nhật ảnh đại diện"
#   - B12 status Missing -> Available
#   - C12 note replaced with new avatar-update note (wrapped)
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value2 = "Available"

$ws.Range("C12").Value2 = "Đã có thể cập nhật ảnh đại diện ở thông tin user và lúc đăng ký `ntài khoản"
$ws.Range("C12").WrapText = $true

$ws.Rows.Item(12).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 11 : "Thay đổi thông tin cá nhân"
#   - B11 status Missing -> Dev
#   - C11 note replaced with new "chưa đổi được mật khẩu" note (wrapped)
# ---------------------------------------------------------------------------
$ws.Range("A22").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("B14").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value2 = "Dev"

$ws.Range("C11").Value2 = "có thể đổi tên nhưng không lưu tên mới ở adminpage, chưa thây đổi `nđược mật khẩu"
$ws.Range("C11").WrapText = $true

$ws.Rows.Item(11).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 22 : "Tìm kiếm cuộc hội thoại ở trang lịch sử" - row just got shorter
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 28.8

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet view : scroll back to the top and move the active selection
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E12").Select()
